$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 67; existing rows 67:81 shift down to 68:82
$ws.Rows("67:67").Insert()

# Update the sheet dimension-affecting new row with the new weekly record.
# Columns A-C, E-I, N, O, Q, R mirror the record that used to sit in row 67
# (now shifted to row 68); only D (fecha), K, L, M, P (prices) differ.
$ws.Range("A67").Value = 1
$ws.Range("B67").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C67").Value = "Arica y Parinacota"
$ws.Range("D67").Value = 45244
$ws.Range("D67").NumberFormat = $ws.Range("D68").NumberFormat
$ws.Range("E67").Value = 15
$ws.Range("F67").Value = 100112052
$ws.Range("G67").Value = "Albahaca"
$ws.Range("H67").Value = "Sin especificar"
$ws.Range("I67").Value = "Primera"
$ws.Range("J67").Value = 300
$ws.Range("K67").Value = 1000
$ws.Range("L67").Value = 1500
$ws.Range("M67").Value = 1250
$ws.Range("N67").Value = "$/paquete"
$ws.Range("O67").Value = "Región de Arica y Parinacota"
$ws.Range("P67").Value = 1250
$ws.Range("Q67").Value = 1
$ws.Range("R67").Value = "Hortaliza"
